$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated test-range results (columns C:M) for every data row (2-18).
# Same new figures are reused across the four repeated label blocks
# (rows 2-5, 6-9, 10-13, 14-17) and the trailing "Maximum" row (18),
# mirroring the structure of the original sheet.
$rowData = @{
    2 = @{ "C" = -0.2; "D" = 0.34; "E" = 0.21; "F" = 0.13; "G" = -0.23; "H" = -0.24; "I" = 0.05; "J" = 0.2; "K" = -0.02; "L" = 0.24; "M" = 0.35 }
    3 = @{ "C" = -0.13; "D" = 0.3; "E" = 0.22; "F" = 0.08; "G" = -0.25; "H" = -0.26; "I" = 0.04; "J" = 0.17; "K" = -0.04; "L" = 0.23; "M" = 0.3 }
    4 = @{ "C" = -0.13; "D" = 0.3; "E" = 0.22; "F" = 0.08; "G" = -0.25; "H" = -0.27; "I" = 0.04; "J" = 0.15; "K" = -0.04; "L" = 0.23; "M" = 0.28999999999999998 }
    5 = @{ "C" = -0.13; "D" = 0.31; "E" = 0.22; "F" = 0.08; "G" = -0.26; "H" = -0.27; "I" = 0.04; "J" = 0.15; "K" = -0.04; "L" = 0.24; "M" = 0.28999999999999998 }
    6 = @{ "C" = -0.24; "D" = 0.09; "E" = -0.39; "F" = 0.04; "G" = -0.22; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = -0.04; "L" = -0.28000000000000003; "M" = 0.37 }
    7 = @{ "C" = -0.23; "D" = -0.09; "E" = -0.62; "F" = -0.17; "G" = -0.22; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = 0.01; "L" = -0.28000000000000003; "M" = 0.24 }
    8 = @{ "C" = -0.24; "D" = -0.11; "E" = -0.73; "F" = -0.17; "G" = -0.22; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = 0; "L" = -0.28000000000000003; "M" = 0.16 }
    9 = @{ "C" = -0.25; "D" = -0.13; "E" = -0.39; "F" = -0.17; "G" = -0.23; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = 0.01; "L" = -0.28000000000000003; "M" = -0.02 }
    10 = @{ "C" = -0.2; "D" = 0.34; "E" = 0.21; "F" = 0.13; "G" = -0.23; "H" = -0.24; "I" = 0.05; "J" = 0.2; "K" = -0.02; "L" = 0.24; "M" = 0.35 }
    11 = @{ "C" = -0.13; "D" = 0.3; "E" = 0.22; "F" = 0.08; "G" = -0.25; "H" = -0.26; "I" = 0.04; "J" = 0.17; "K" = -0.04; "L" = 0.23; "M" = 0.3 }
    12 = @{ "C" = -0.13; "D" = 0.3; "E" = 0.22; "F" = 0.08; "G" = -0.25; "H" = -0.27; "I" = 0.04; "J" = 0.15; "K" = -0.04; "L" = 0.23; "M" = 0.28999999999999998 }
    13 = @{ "C" = -0.13; "D" = 0.31; "E" = 0.22; "F" = 0.08; "G" = -0.26; "H" = -0.27; "I" = 0.04; "J" = 0.15; "K" = -0.04; "L" = 0.24; "M" = 0.28999999999999998 }
    14 = @{ "C" = -0.24; "D" = 0.09; "E" = -0.39; "F" = 0.04; "G" = -0.22; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = -0.04; "L" = -0.28000000000000003; "M" = 0.37 }
    15 = @{ "C" = -0.23; "D" = -0.09; "E" = -0.62; "F" = -0.17; "G" = -0.22; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = 0.01; "L" = -0.28000000000000003; "M" = 0.24 }
    16 = @{ "C" = -0.24; "D" = -0.11; "E" = -0.73; "F" = -0.17; "G" = -0.22; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = 0; "L" = -0.28000000000000003; "M" = 0.16 }
    17 = @{ "C" = -0.25; "D" = -0.13; "E" = -0.39; "F" = -0.17; "G" = -0.23; "H" = -0.23; "I" = 0.03; "J" = -0.05; "K" = 0.01; "L" = -0.28000000000000003; "M" = -0.02 }
    18 = @{ "C" = -0.13; "D" = 0.34; "E" = 0.22; "F" = 0.13; "G" = -0.22; "H" = -0.23; "I" = 0.05; "J" = 0.2; "K" = 0.01; "L" = 0.24; "M" = 0.37 }
}

foreach ($rowNum in $rowData.Keys) {
    $cols = $rowData[$rowNum]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$rowNum").Value = $cols[$col]
    }
}
